$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5098
$ws.Range("J40").Value = 6285.7144
$ws.Range("L40").Value = 6285.7144
$ws.Range("N40").Value = -6635.7144

$ws.Range("H74").Value = 6741.0586
$ws.Range("I74").Value = 5312.857
$ws.Range("J74").Value = 7740.8
$ws.Range("K74").Value = 5312.857
$ws.Range("L74").Value = 7740.8
$ws.Range("M74").Value = -4376.857
$ws.Range("N74").Value = -9612.799999999999

$ws.Range("H77").Value = 6741.0586
$ws.Range("I77").Value = 5312.857
$ws.Range("J77").Value = 7740.8
$ws.Range("K77").Value = 26564.285
$ws.Range("L77").Value = 38704
$ws.Range("M77").Value = -21884.285
$ws.Range("N77").Value = -48064

$ws.Range("H93").Value = 40890
$ws.Range("J93").Value = 40890
$ws.Range("L93").Value = 40890
$ws.Range("N93").Value = -45882

$ws.Range("H112").Value = 1336.1628
$ws.Range("J112").Value = 1348.9286
$ws.Range("L112").Value = 4046.7858
$ws.Range("N112").Value = -6262.7858

$ws.Range("H129").Value = 1260.0588
$ws.Range("J129").Value = 1320.5696
$ws.Range("L129").Value = 3961.7088
$ws.Range("N129").Value = -13961.7088

$ws.Range("H137").Value = 605433.7
$ws.Range("I137").Value = 1109787.1
$ws.Range("J137").Value = 3011.5833
$ws.Range("K137").Value = 3329361.3
$ws.Range("L137").Value = 9034.749899999999
$ws.Range("M137").Value = -3326811.3
$ws.Range("N137").Value = -14134.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 207.6
$ws.Range("I5").Value = 135
$ws.Range("J5").Value = 498
$ws.Range("K5").Value = 135
$ws.Range("L5").Value = 498
$ws.Range("M5").Value = -23
$ws.Range("N5").Value = -722

$ws.Range("H61").Value = 1809.0769
$ws.Range("I61").Value = 1851.5
$ws.Range("J61").Value = 1300
$ws.Range("K61").Value = 1851.5
$ws.Range("L61").Value = 1300
$ws.Range("M61").Value = -1639.5
$ws.Range("N61").Value = -1724

$ws.Range("H63").Value = 9238271
$ws.Range("I63").Value = 19789466
$ws.Range("J63").Value = 5975
$ws.Range("K63").Value = 19789466
$ws.Range("L63").Value = 5975
$ws.Range("M63").Value = -19788780
$ws.Range("N63").Value = -7347

$ws.Range("H66").Value = 9238271
$ws.Range("I66").Value = 19789466
$ws.Range("J66").Value = 5975
$ws.Range("K66").Value = 98947330
$ws.Range("L66").Value = 29875
$ws.Range("M66").Value = -98943898
$ws.Range("N66").Value = -36739

$ws.Range("H74").Value = 218627.8
$ws.Range("I74").Value = 368053.66
$ws.Range("J74").Value = 1960.3
$ws.Range("K74").Value = 368053.66
$ws.Range("L74").Value = 1960.3
$ws.Range("M74").Value = -367179.66
$ws.Range("N74").Value = -3708.3

$ws.Range("H77").Value = 218627.8
$ws.Range("I77").Value = 368053.66
$ws.Range("J77").Value = 1960.3
$ws.Range("K77").Value = 1840268.3
$ws.Range("L77").Value = 9801.5
$ws.Range("M77").Value = -1835900.3
$ws.Range("N77").Value = -18537.5

$ws.Range("H102").Value = 1436.6666
$ws.Range("I102").Value = 1436.6666
$ws.Range("K102").Value = 1436.6666
$ws.Range("M102").Value = 185.3334

$ws.Range("H103").Value = 33444.445
$ws.Range("J103").Value = 33444.445
$ws.Range("L103").Value = 33444.445
$ws.Range("N103").Value = -35788.445

$ws.Range("H136").Value = 1809.0769
$ws.Range("I136").Value = 1851.5
$ws.Range("J136").Value = 1300
$ws.Range("K136").Value = 5554.5
$ws.Range("L136").Value = 3900
$ws.Range("M136").Value = -3004.5
$ws.Range("N136").Value = -9000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 207.6
$ws.Range("I4").Value = 135
$ws.Range("J4").Value = 498
$ws.Range("K4").Value = 135
$ws.Range("L4").Value = 498
$ws.Range("M4").Value = -20
$ws.Range("N4").Value = -728

$ws.Range("H95").Value = 32625
$ws.Range("J95").Value = 32625
$ws.Range("L95").Value = 32625
$ws.Range("N95").Value = -38117

$ws.Range("M97").Value = $null
$ws.Range("H97").Value = 4500
$ws.Range("I97").Value = 4500
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 4500
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = -3509

$ws.Range("H103").Value = 34700
$ws.Range("J103").Value = 34700
$ws.Range("L103").Value = 34700
$ws.Range("N103").Value = -37044

$ws.Range("H134").Value = 2521.712
$ws.Range("I134").Value = 853.7179599999999
$ws.Range("J134").Value = 5774.3
$ws.Range("K134").Value = 2561.15388
$ws.Range("L134").Value = 17322.9
$ws.Range("M134").Value = -26.15387999999984
$ws.Range("N134").Value = -22392.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 189666.12
$ws.Range("I31").Value = 436565.6
$ws.Range("J31").Value = 2986.0488
$ws.Range("K31").Value = 436565.6
$ws.Range("L31").Value = 2986.0488
$ws.Range("M31").Value = -436270.6
$ws.Range("N31").Value = -3576.0488

$ws.Range("H34").Value = 189666.12
$ws.Range("I34").Value = 436565.6
$ws.Range("J34").Value = 2986.0488
$ws.Range("K34").Value = 436565.6
$ws.Range("L34").Value = 2986.0488
$ws.Range("M34").Value = -436363.6
$ws.Range("N34").Value = -3390.0488

$ws.Range("H58").Value = 2898.3447
$ws.Range("I58").Value = 1502.2609
$ws.Range("K58").Value = 1502.2609
$ws.Range("M58").Value = -1299.2609

$ws.Range("H99").Value = 4689.8667
$ws.Range("J99").Value = 6508
$ws.Range("L99").Value = 6508
$ws.Range("N99").Value = -9504

$ws.Range("H107").Value = 2703359.5
$ws.Range("I107").Value = 4000333.5
$ws.Range("K107").Value = 4000333.5
$ws.Range("M107").Value = -3998413.5

$ws.Range("H126").Value = 4689.8667
$ws.Range("J126").Value = 6508
$ws.Range("L126").Value = 19524
$ws.Range("N126").Value = -24464

$ws.Range("H132").Value = 3514.724
$ws.Range("I132").Value = 2928.5
$ws.Range("J132").Value = 5357.143
$ws.Range("K132").Value = 8785.5
$ws.Range("L132").Value = 16071.429
$ws.Range("M132").Value = -6255.5
$ws.Range("N132").Value = -21131.429

$ws.Range("H136").Value = 2898.3447
$ws.Range("I136").Value = 1502.2609
$ws.Range("K136").Value = 4506.7827
$ws.Range("M136").Value = -1956.7827

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3902.1667
$ws.Range("I68").Value = 1341.091
$ws.Range("J68").Value = 4810.9355
$ws.Range("K68").Value = 4023.273
$ws.Range("L68").Value = 14432.8065
$ws.Range("M68").Value = -3212.273
$ws.Range("N68").Value = -16054.8065

$ws.Range("H71").Value = 3902.1667
$ws.Range("I71").Value = 1341.091
$ws.Range("J71").Value = 4810.9355
$ws.Range("K71").Value = 12069.819
$ws.Range("L71").Value = 43298.4195
$ws.Range("M71").Value = -8013.819
$ws.Range("N71").Value = -51410.4195

$ws.Range("H137").Value = 2235.6316
$ws.Range("J137").Value = 2037.5555
$ws.Range("L137").Value = 6112.666499999999
$ws.Range("N137").Value = -16312.6665

$ws.Range("H140").Value = 2424.348
$ws.Range("I140").Value = 695.55554
$ws.Range("J140").Value = 3535.7144
$ws.Range("K140").Value = 2086.66662
$ws.Range("L140").Value = 10607.1432
$ws.Range("M140").Value = 3093.33338
$ws.Range("N140").Value = -20967.1432

$ws.Range("H141").Value = 4882.2856
$ws.Range("I141").Value = 4729.3335
$ws.Range("K141").Value = 14188.0005
$ws.Range("M141").Value = -9008.000499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 14785.105
$ws.Range("I43").Value = 1531.4
$ws.Range("J43").Value = 29511.445
$ws.Range("K43").Value = 1531.4
$ws.Range("L43").Value = 29511.445
$ws.Range("M43").Value = -1380.4
$ws.Range("N43").Value = -29813.445

$ws.Range("H46").Value = 23487.268
$ws.Range("J46").Value = 23361.357
$ws.Range("L46").Value = 23361.357
$ws.Range("N46").Value = -23673.357

$ws.Range("H57").Value = 32246.666
$ws.Range("J57").Value = 32246.666
$ws.Range("L57").Value = 32246.666
$ws.Range("N57").Value = -33886.666

$ws.Range("H80").Value = 2795
$ws.Range("I80").Value = 2702.7273
$ws.Range("J80").Value = 3133.3333
$ws.Range("K80").Value = 2702.7273
$ws.Range("L80").Value = 3133.3333
$ws.Range("M80").Value = -1704.7273
$ws.Range("N80").Value = -5129.3333

$ws.Range("H83").Value = 2795
$ws.Range("I83").Value = 2702.7273
$ws.Range("J83").Value = 3133.3333
$ws.Range("K83").Value = 13513.6365
$ws.Range("L83").Value = 15666.6665
$ws.Range("M83").Value = -8521.636500000001
$ws.Range("N83").Value = -25650.6665

$ws.Range("H113").Value = 1695.25
$ws.Range("I113").Value = 1695.25
$ws.Range("K113").Value = 1695.25
$ws.Range("M113").Value = 474.75

$ws.Range("H122").Value = 6066.9
$ws.Range("I122").Value = 5124.875
$ws.Range("J122").Value = 9835
$ws.Range("K122").Value = 15374.625
$ws.Range("L122").Value = 29505
$ws.Range("M122").Value = -12924.625
$ws.Range("N122").Value = -34405

$ws.Range("H132").Value = 3274.1292
$ws.Range("I132").Value = 2490
$ws.Range("J132").Value = 5528.5
$ws.Range("K132").Value = 7470
$ws.Range("L132").Value = 16585.5
$ws.Range("M132").Value = -4940
$ws.Range("N132").Value = -21645.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4092.0881
$ws.Range("I132").Value = 3006.1052
$ws.Range("J132").Value = 5467.6665
$ws.Range("K132").Value = 9018.3156
$ws.Range("L132").Value = 16402.9995
$ws.Range("M132").Value = -6488.3156
$ws.Range("N132").Value = -21462.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3186.8572
$ws.Range("I122").Value = 1768.7222
$ws.Range("J122").Value = 4688.4116
$ws.Range("K122").Value = 5306.1666
$ws.Range("L122").Value = 14065.2348
$ws.Range("M122").Value = -2856.1666
$ws.Range("N122").Value = -18965.2348

$ws.Range("H132").Value = 1717.2727
$ws.Range("I132").Value = 732.2963
$ws.Range("K132").Value = 2196.8889
$ws.Range("M132").Value = 333.1111000000001

$ws.Range("H136").Value = 2283.4814
$ws.Range("I136").Value = 940.25
$ws.Range("K136").Value = 2820.75
$ws.Range("M136").Value = -270.75
